# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-27 from
# serial date 45315 (2024-01-24) to 45316 (2024-01-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45315) {
        $cell.Value2 = 45316
    }
}
